# Water level historical values and norm data support changes
#
# The "discharge" template workbook gets a new "water_level" worksheet
# added right after the existing "discharge" sheet. The new sheet is a
# duplicate of "discharge" (same header row "Period", decadal period
# numbers 1-36, and "Value" label) and becomes the active/selected sheet.

$wb = $excel.ActiveWorkbook

# Locate the existing "discharge" worksheet.
$discharge = $wb.Worksheets.Item("discharge")

# Duplicate it immediately after itself - this carries over the header
# row, the "Value" cell, column/row formatting and cell styles.
$discharge.Copy($null, $discharge)

# The copy is inserted right after "discharge" and is given a default
# name like "discharge (2)" - rename it to "water_level".
$waterLevel = $wb.Worksheets.Item(2)
$waterLevel.Name = "water_level"

# Normalize row 1 height back to the sheet default (15pt) on both sheets
# instead of the inherited custom 13.8pt height.
$discharge.Rows.Item(1).AutoFit()
$waterLevel.Rows.Item(1).AutoFit()

# Make "water_level" the active sheet / tab (matches activeTab="1").
$waterLevel.Activate()
$waterLevel.Select()
$waterLevel.Range("B2").Select()
